$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data for a Nissan Juke (no Year/HP provided)
$ws.Range("A5").Value = "Nissan"
$ws.Range("B5").Value = "Juke"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "B"

# Update the active selection to match the edited workbook state
$ws.Range("F8").Select()
